# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns for rows 2-48,
# and fully rewrites rows 49-51 (Coin/Link/Price/Volume) to reflect
# a newly inserted "BabyDogeCoin" row pushing Algorand/Cronos down
# and dropping the trailing "USDD" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text (avoids Excel auto-converting
# numeric-looking strings like "208.00" into the number 208), while
# resetting the cell style afterwards so no stray "Text"/quote-prefix
# style gets attached to the cell.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# --- Update Price / Volume(1h) for existing rows ---
Set-TextValue $ws.Cells.Item(2, 4) "27.490.85"
$ws.Cells.Item(2, 5).Value = "  -0.10%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.569.54"
$ws.Cells.Item(3, 5).Value = "  -0.55%  "
$ws.Cells.Item(4, 5).Value = "  -0.44%  "
Set-TextValue $ws.Cells.Item(5, 4) "208.00"
$ws.Cells.Item(5, 5).Value = "  +0.91%  "
$ws.Cells.Item(6, 5).Value = "  -0.76%  "
$ws.Cells.Item(7, 5).Value = "  -0.38%  "
Set-TextValue $ws.Cells.Item(8, 4) "22.02"
$ws.Cells.Item(8, 5).Value = "  -0.47%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.249"
$ws.Cells.Item(9, 5).Value = "  -1.16%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.0590"
$ws.Cells.Item(10, 5).Value = "  +0.21%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0867"
$ws.Cells.Item(11, 5).Value = "  +0.23%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.792.48"
$ws.Cells.Item(12, 5).Value = "  -0.70%  "
Set-TextValue $ws.Cells.Item(13, 4) "1.573.50"
$ws.Cells.Item(13, 5).Value = "  -0.38%  "
$ws.Cells.Item(14, 5).Value = "  -0.27%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.519"
$ws.Cells.Item(15, 5).Value = "  -2.28%  "
Set-TextValue $ws.Cells.Item(16, 4) "63.27"
$ws.Cells.Item(16, 5).Value = "  +0.59%  "
Set-TextValue $ws.Cells.Item(17, 4) "27.477.29"
$ws.Cells.Item(17, 5).Value = "  -0.21%  "
Set-TextValue $ws.Cells.Item(18, 4) "214.12"
$ws.Cells.Item(18, 5).Value = "  -0.91%  "
Set-TextValue $ws.Cells.Item(19, 4) "0.0₃0690"
$ws.Cells.Item(19, 5).Value = "  +0.11%  "
Set-TextValue $ws.Cells.Item(20, 4) "7.26"
$ws.Cells.Item(20, 5).Value = "  -0.71%  "
$ws.Cells.Item(21, 5).Value = "  -0.43%  "
$ws.Cells.Item(22, 5).Value = "  -0.27%  "
Set-TextValue $ws.Cells.Item(23, 4) "9.56"
$ws.Cells.Item(23, 5).Value = "  +0.65%  "
Set-TextValue $ws.Cells.Item(24, 4) "2.02"
$ws.Cells.Item(24, 5).Value = "  +1.29%  "
Set-TextValue $ws.Cells.Item(25, 4) "153.13"
$ws.Cells.Item(25, 5).Value = "  -0.14%  "
$ws.Cells.Item(26, 5).Value = "  +1.70%  "
$ws.Cells.Item(27, 5).Value = "  -0.33%  "
Set-TextValue $ws.Cells.Item(28, 4) "15.02"
$ws.Cells.Item(28, 5).Value = "  +0.04%  "
$ws.Cells.Item(29, 5).Value = "  -1.21%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.16"
$ws.Cells.Item(30, 5).Value = "  -0.13%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.0471"
$ws.Cells.Item(31, 5).Value = "  +1.21%  "
Set-TextValue $ws.Cells.Item(32, 4) "3.20"
$ws.Cells.Item(32, 5).Value = "  -1.03%  "
Set-TextValue $ws.Cells.Item(33, 4) "1.364.32"
$ws.Cells.Item(33, 5).Value = "  +0.06%  "
Set-TextValue $ws.Cells.Item(34, 4) "2.95"
$ws.Cells.Item(34, 5).Value = "  +0.39%  "
$ws.Cells.Item(35, 5).Value = "  +2.62%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.976"
$ws.Cells.Item(36, 5).Value = "  +0.88%  "
$ws.Cells.Item(37, 5).Value = "  -0.21%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.0168"
$ws.Cells.Item(38, 5).Value = "  +1.87%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.533"
$ws.Cells.Item(39, 5).Value = "  -0.43%  "
Set-TextValue $ws.Cells.Item(40, 4) "0.822"
$ws.Cells.Item(40, 5).Value = "  +1.82%  "
$ws.Cells.Item(41, 5).Value = "  -0.36%  "
$ws.Cells.Item(42, 5).Value = "  -0.08%  "
$ws.Cells.Item(43, 5).Value = "  +1.81%  "
Set-TextValue $ws.Cells.Item(44, 4) "64.17"
$ws.Cells.Item(44, 5).Value = "  +1.37%  "
Set-TextValue $ws.Cells.Item(45, 4) "5.28"
$ws.Cells.Item(45, 5).Value = "  +0.29%  "
$ws.Cells.Item(46, 5).Value = "  -1.16%  "
Set-TextValue $ws.Cells.Item(47, 4) "1.706.05"
$ws.Cells.Item(47, 5).Value = "  -0.79%  "
Set-TextValue $ws.Cells.Item(48, 4) "85.47"
$ws.Cells.Item(48, 5).Value = "  -1.58%  "

# --- Rows 49-51: BabyDogeCoin inserted, Algorand/Cronos shifted down, USDD dropped ---
$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Cells.Item(49, 4) "0.0₇0987"
$ws.Cells.Item(49, 5).Value = "  +2.49%  "

$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Cells.Item(50, 4) "0.0955"
$ws.Cells.Item(50, 5).Value = "  -1.12%  "

$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(51, 4) "0.0495"
$ws.Cells.Item(51, 5).Value = "  -0.14%  "

